# The sheet's "Docentes responsáveis:" section onward was mis-aligned: each
# label in column A was sitting one row above the row that actually held its
# value (so e.g. "Objetivos:" displayed the docente's name, "Método:" showed
# the docente's name again, "Critério:" showed the evaluation method, etc.),
# and the "Bibliografia:" field had no value at all.
#
# Fix this by inserting one new row right after the "Docentes responsáveis:"
# label (current row 13), which shifts the mis-aligned label/value rows back
# into register, then fill in the correct value for every field that was
# wrong (or previously missing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at row 13 (shifts rows 13-24 down to 14-25).
$ws.Rows.Item(13).Insert()

# The inserted row copies formatting from the row above, which leaves a
# stray, empty, styled cell in A13 and gives B13/C13 the wrong style. Copy
# the correct data-cell formatting (from the row that used to be 13, now
# shifted to row 14) onto B13:C13, then drop the stray A13 cell entirely.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$excel.CutCopyMode = 0

# 2) "Objetivos:" (row 10) incorrectly held the docente's name; set the
#    real Portuguese objectives text.
$ws.Range("B10").Value = 'Verificação experimental dos conceitos básicos de hidrostática, hidrodinâmica, termodinâmica e ondas.'
$ws.Range("C10").Value = 'Verificação experimental dos conceitos básicos de hidrostática, hidrodinâmica, termodinâmica e ondas.'

# 3) Fill the newly inserted row 13 (under "Docentes responsáveis:") with
#    the responsible docente's name.
$ws.Range("B13").Value = '5817535 - Lucas Barboza Sarno da Silva'
$ws.Range("C13").Value = '5817535 - Lucas Barboza Sarno da Silva'

# 4) "Programa resumido:" (row 14) incorrectly held "Semestral"; set the
#    real Portuguese short-syllabus text.
$ws.Range("B14").Value = 'Abordagem experimental de conceitos relacionados à mecânica dos fluidos, termodinâmica, oscilações e ondas.'
$ws.Range("C14").Value = 'Abordagem experimental de conceitos relacionados à mecânica dos fluidos, termodinâmica, oscilações e ondas.'

# 5) "Programa:" (row 16) incorrectly held a date; set the real Portuguese
#    syllabus text.
$ws.Range("B16").Value = '1) Princípio de Stevin e Pascal2) Empuxo e Princípio de Arquimedes3) Tensão superficial4) Queda em um meio viscoso5) Sistema massa-mola6) Ondas mecânicas7) Calor, temperatura e capacidade do corpo de armazenar energia8) Dilatação linear9) Os meios de propagação de calor10) Calor específico e calor latente11) A lei de Boyle-Mariotte'
$ws.Range("C16").Value = '1) Princípio de Stevin e Pascal2) Empuxo e Princípio de Arquimedes3) Tensão superficial4) Queda em um meio viscoso5) Sistema massa-mola6) Ondas mecânicas7) Calor, temperatura e capacidade do corpo de armazenar energia8) Dilatação linear9) Os meios de propagação de calor10) Calor específico e calor latente11) A lei de Boyle-Mariotte'

# 6) "Método:" (row 19) incorrectly held the docente's name again; set the
#    real evaluation-method text.
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# 7) "Critério:" (row 20) incorrectly held the evaluation-method text; set
#    the real passing-criterion text.
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'

# 8) "Norma de recuperação:" (row 21) incorrectly held the passing-criterion
#    text; set the real recovery-norm text.
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'

# 9) "Bibliografia:" (row 22) incorrectly held the recovery-norm text; set
#    the real bibliography text.
$ws.Range("B22").Value = '1. Apostilas do Laboratório de Ensino de Física do IFSC/USP.2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).3. NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).4. RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 2, LTC (2008).5. TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 2, LTC (2008).6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física II, Vol. 2,     Pearson Addison Wesley (2009).7. JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008)'
$ws.Range("C22").Value = '1. Apostilas do Laboratório de Ensino de Física do IFSC/USP.2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).3. NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).4. RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 2, LTC (2008).5. TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 2, LTC (2008).6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física II, Vol. 2,     Pearson Addison Wesley (2009).7. JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008)'
